# Insert a new data record at row 23 (pushing the existing rows 23:124 down
# to 24:125, which is how the underlying diff shows every subsequent row's
# values shifting down by exactly one row, with a brand-new final row 125
# appearing that matches the old row 124).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

$ws.Cells.Item(23, 1).Value = 1
$ws.Cells.Item(23, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(23, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(23, 4).Value = 44859
$ws.Cells.Item(23, 5).Value = 15
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100102
$ws.Cells.Item(23, 8).Value = "Cítricos"
$ws.Cells.Item(23, 9).Value = 100102004
$ws.Cells.Item(23, 10).Value = "Mandarina"
$ws.Cells.Item(23, 11).Value = "Murcott"
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 250
$ws.Cells.Item(23, 14).Value = 15000
$ws.Cells.Item(23, 15).Value = 16000
$ws.Cells.Item(23, 16).Value = 15500
$ws.Cells.Item(23, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(23, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(23, 19).Value = 775
$ws.Cells.Item(23, 20).Value = 20
